# Auto-generated Excel COM-interop edit script
# Applies the Go-Live sync update described in the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 00_읽는법 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = 'S1(20260221 하드닝), S2(202603XX 종합보고서), S3(Go-Live Gap Closure), S4(spec_sync_report), S5(golive 증적)'
$ws1.Range("A6").Value = '최신 반영 범위'
$ws1.Range("B6").Value = 'Node 22 표준화, role/access_level 정규화, MFA/세션, RBAC 2인 승인, audit chain verify, runbook 3종'

# --- Sheet 2: 01_용어사전 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A12").Value = 'access_level'
$ws2.Range("B12").Value = '접근수준 표기'
$ws2.Range("C12").Value = 'ROLE이 아닌 공개/인증필요 접근 속성'
$ws2.Range("D12").Value = 'ROLE 혼용 방지'
$ws2.Range("E12").Value = 'API 스펙 비고'
$ws2.Range("F12").Value = '권한 칼럼과 구분해 읽기'
$ws2.Range("G12").Value = '권한(ROLE)과 혼동 금지'
$ws2.Range("H12").Value = '-'
$ws2.Range("A13").Value = 'MFA'
$ws2.Range("B13").Value = '다중 인증'
$ws2.Range("C13").Value = '비밀번호 외에 추가 인증(TOTP 등)을 요구하는 보안 방식'
$ws2.Range("D13").Value = '관리자 계정 탈취 위험 감소'
$ws2.Range("E13").Value = 'OPS/ADMIN 로그인'
$ws2.Range("F13").Value = '코드 입력/복구코드 보관'
$ws2.Range("G13").Value = '반복 실패 시 잠금 가능'
$ws2.Range("H13").Value = 'AUTH_MFA_INVALID_CODE/AUTH_MFA_LOCKED'
$ws2.Range("A14").Value = 'audit chain verify'
$ws2.Range("B14").Value = '감사체인 무결성 점검'
$ws2.Range("C14").Value = '감사로그 hash chain이 끊기지 않았는지 검증하는 절차/API'
$ws2.Range("D14").Value = '위변조 조기 탐지'
$ws2.Range("E14").Value = '운영 점검/감사'
$ws2.Range("F14").Value = '이상 시 runbook 따라 대응'
$ws2.Range("G14").Value = '검증 실패 시 export 제한 검토'
$ws2.Range("H14").Value = '-'
$ws2.Range("A15").Value = 'session revoke'
$ws2.Range("B15").Value = '세션 강제 종료'
$ws2.Range("C15").Value = '내 세션 목록에서 의심 세션을 종료하거나 타 세션을 일괄 종료'
$ws2.Range("D15").Value = '계정 도난 대응 속도 향상'
$ws2.Range("E15").Value = '보안 설정 화면'
$ws2.Range("F15").Value = '의심 세션 즉시 종료'
$ws2.Range("G15").Value = '현재 세션과 타 세션 구분'
$ws2.Range("H15").Value = '-'
$ws2.Range("A16").Value = 'runbook'
$ws2.Range("B16").Value = '운영 대응 절차서'
$ws2.Range("C16").Value = '장애/보안 이슈 발생 시 단계별 확인·복구 절차를 정리한 문서'
$ws2.Range("D16").Value = '사고 대응 표준화'
$ws2.Range("E16").Value = '운영 문서'
$ws2.Range("F16").Value = '장애시 문서 순서대로 수행'
$ws2.Range("G16").Value = '임의 조치 금지'
$ws2.Range("H16").Value = '-'

# --- Sheet 3: 02_역할권한표 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C4").Value = 'RBAC 변경요청 생성/승인'
$ws3.Range("D4").Value = 'OPS 전용 차단조치 직접 실행'
$ws3.Range("E4").Value = 'RBAC Matrix/Approval'
$ws3.Range("G4").Value = '직접 적용이 아닌 승인흐름 준수'
$ws3.Range("H4").Value = 'S1,S2,S3'
$ws3.Range("C5").Value = '지표/감사/차단/체인검증'
$ws3.Range("D5").Value = 'RBAC 정책 승인'
$ws3.Range("E5").Value = 'Dashboard/Audit/Chain Verify'
$ws3.Range("G5").Value = 'export 범위 제한 준수'
$ws3.Range("H5").Value = 'S1,S2,S3'

# --- Sheet 4: 03_메뉴_기능맵 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C6").Value = '권한 변경 요청 생성'
$ws4.Range("E6").Value = '변경요청(PENDING)'
$ws4.Range("F6").Value = 'SEC-002-403'
$ws4.Range("G6").Value = '직접 적용이 아닌 승인 대기'
$ws4.Range("H6").Value = 'S2,S3'
$ws4.Range("A8").Value = 'RBAC 승인요청'
$ws4.Range("B8").Value = 'ADMIN(SYSTEM_ADMIN)'
$ws4.Range("C8").Value = '요청 승인/반려'
$ws4.Range("D8").Value = 'request_id/사유'
$ws4.Range("E8").Value = 'APPROVED/REJECTED'
$ws4.Range("F8").Value = 'SEC-002-403'
$ws4.Range("G8").Value = '요청자 본인승인 금지'
$ws4.Range("H8").Value = 'S2,S3'
$ws4.Range("A9").Value = 'MFA 설정'
$ws4.Range("B9").Value = 'OPS,ADMIN'
$ws4.Range("C9").Value = 'TOTP 등록/활성화'
$ws4.Range("D9").Value = '인증앱 코드'
$ws4.Range("E9").Value = 'MFA 활성 상태'
$ws4.Range("F9").Value = 'AUTH_MFA_SETUP_REQUIRED'
$ws4.Range("G9").Value = '복구코드 안전보관'
$ws4.Range("H9").Value = 'S2'
$ws4.Range("A10").Value = '내 세션 관리'
$ws4.Range("B10").Value = '인증 사용자'
$ws4.Range("C10").Value = '세션 목록 조회/강제종료'
$ws4.Range("D10").Value = '세션 선택'
$ws4.Range("E10").Value = 'revoke 결과'
$ws4.Range("F10").Value = 'SEC-001-401'
$ws4.Range("G10").Value = '의심 세션 즉시 revoke'
$ws4.Range("H10").Value = 'S2'
$ws4.Range("A11").Value = 'Audit Chain Verify'
$ws4.Range("B11").Value = 'OPS'
$ws4.Range("C11").Value = '감사 체인 무결성 점검'
$ws4.Range("D11").Value = 'tenant/기간'
$ws4.Range("E11").Value = '검증 PASS/FAIL'
$ws4.Range("F11").Value = 'SEC-002-403'
$ws4.Range("G11").Value = '실패 시 runbook 즉시 실행'
$ws4.Range("H11").Value = 'S2'

# --- Sheet 5: 04_동작흐름_한눈표 ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A6").Value = '''5'
$ws5.Range("B6").Value = 'OPS/ADMIN 로그인'
$ws5.Range("C6").Value = 'MFA 필요 여부/챌린지 검증'
$ws5.Range("D6").Value = '2차 인증 완료 후 로그인'
$ws5.Range("E6").Value = '401/429'
$ws5.Range("F6").Value = 'AUTH_MFA_INVALID_CODE/AUTH_MFA_LOCKED'
$ws5.Range("G6").Value = 'S2'
$ws5.Range("A7").Value = '''6'
$ws5.Range("B7").Value = 'RBAC 변경 요청'
$ws5.Range("C7").Value = '요청 생성 -> 2인 승인'
$ws5.Range("D7").Value = '승인 후 정책 반영'
$ws5.Range("E7").Value = '403/409'
$ws5.Range("F7").Value = 'SEC-002-403'
$ws5.Range("G7").Value = 'S2,S3'
$ws5.Range("A8").Value = '''7'
$ws5.Range("B8").Value = '감사 무결성 점검'
$ws5.Range("C8").Value = 'chain-verify 실행'
$ws5.Range("D8").Value = 'PASS 보고'
$ws5.Range("E8").Value = 'FAIL 경보/조치'
$ws5.Range("F8").Value = '-'
$ws5.Range("G8").Value = 'S2'

# --- Sheet 6: 05_오류코드_쉬운설명 ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A7").Value = 'AUTH_MFA_SETUP_REQUIRED'
$ws6.Range("B7").Value = 'MFA 설정 필요'
$ws6.Range("C7").Value = 'OPS/ADMIN 계정에 MFA 미설정'
$ws6.Range("D7").Value = 'MFA 등록 진행'
$ws6.Range("E7").Value = '계정 보안정책 확인'
$ws6.Range("F7").Value = '-'
$ws6.Range("G7").Value = '관리자 MFA 필수'
$ws6.Range("H7").Value = 'S2'
$ws6.Range("A8").Value = 'AUTH_MFA_INVALID_CODE'
$ws6.Range("B8").Value = 'MFA 코드 불일치'
$ws6.Range("C8").Value = 'TOTP 또는 복구코드 오류'
$ws6.Range("D8").Value = '코드 재확인'
$ws6.Range("E8").Value = '시간 동기화 확인'
$ws6.Range("F8").Value = '-'
$ws6.Range("G8").Value = 'MFA 검증'
$ws6.Range("H8").Value = 'S2'
$ws6.Range("A9").Value = 'AUTH_MFA_LOCKED'
$ws6.Range("B9").Value = 'MFA 입력 잠금'
$ws6.Range("C9").Value = '반복 실패 누적'
$ws6.Range("D9").Value = 'Retry-After 후 재시도'
$ws6.Range("E9").Value = '실패 패턴 점검'
$ws6.Range("F9").Value = '잠금 해제까지'
$ws6.Range("G9").Value = 'MFA brute-force 방어'
$ws6.Range("H9").Value = 'S2'

# --- Sheet 7: 06_주의사항_체크리스트 ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A6").Value = 'Node 22 표준 환경 확인'
$ws7.Range("B6").Value = '로컬/CI 재현성 확보'
$ws7.Range("C6").Value = '.nvmrc 기준 버전 확인'
$ws7.Range("D6").Value = '미완료'
$ws7.Range("E6").Value = 'S2,S3'
$ws7.Range("A7").Value = '권한/접근수준 구분 확인'
$ws7.Range("B7").Value = 'ROLE 혼동으로 인한 오판 방지'
$ws7.Range("C7").Value = '권한=ROLE, 비고=access_level 확인'
$ws7.Range("D7").Value = '미완료'
$ws7.Range("E7").Value = 'S3,S4'
$ws7.Range("A8").Value = 'Runbook 최신본 확인'
$ws7.Range("B8").Value = '장애 대응 속도/정확도 확보'
$ws7.Range("C8").Value = 'scheduler/audit/spec-notion runbook 숙지'
$ws7.Range("D8").Value = '미완료'
$ws7.Range("E8").Value = 'S2'
$ws7.Range("A9").Value = 'Notion 동기화 증적 확인'
$ws7.Range("B9").Value = '스펙-문서 분기 방지'
$ws7.Range("C9").Value = 'spec_sync_report와 상태파일 점검'
$ws7.Range("D9").Value = '미완료'
$ws7.Range("E9").Value = 'S4,S6'

# --- Sheet 8: 07_FAQ_빠른답변 ---
$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("A6").Value = '왜 관리자 로그인에서 인증을 두 번 하나요?'
$ws8.Range("B6").Value = 'MFA 보안정책 때문입니다.'
$ws8.Range("C6").Value = 'OPS/ADMIN 계정은 탈취 위험이 높아 비밀번호 외 2차 인증을 필수로 적용합니다.'
$ws8.Range("D6").Value = 'S2'
$ws8.Range("A7").Value = '권한 변경이 바로 적용되지 않는 이유는?'
$ws8.Range("B7").Value = '2인 승인 정책 때문입니다.'
$ws8.Range("C7").Value = 'RBAC 변경은 요청 생성 후 서로 다른 승인자 2명이 승인해야 실제 반영됩니다.'
$ws8.Range("D7").Value = 'S2,S3'
$ws8.Range("A8").Value = 'Audit Chain Verify는 언제 쓰나요?'
$ws8.Range("B8").Value = '감사로그 무결성 점검 시 사용합니다.'
$ws8.Range("C8").Value = '정기 점검이나 이상 징후 발생 시 체인 검증을 실행하고 실패 시 runbook 절차로 대응합니다.'
$ws8.Range("D8").Value = 'S2'

# --- Sheet 9: 08_문제해결_증상별 ---
$ws9 = $wb.Worksheets.Item(9)
$ws9.Range("A6").Value = 'MFA 코드가 계속 틀렸다고 나옴'
$ws9.Range("B6").Value = '코드 불일치/시간 오차'
$ws9.Range("C6").Value = 'AUTH_MFA_INVALID_CODE'
$ws9.Range("D6").Value = '기기 시간 자동동기화 후 재시도'
$ws9.Range("E6").Value = 'AUTH_MFA_INVALID_CODE'
$ws9.Range("F6").Value = 'S2'
$ws9.Range("A7").Value = 'MFA가 잠겨 로그인 불가'
$ws9.Range("B7").Value = '반복 실패로 잠금'
$ws9.Range("C7").Value = 'Retry-After 확인'
$ws9.Range("D7").Value = '잠금 해제 후 재시도/운영자 확인'
$ws9.Range("E7").Value = 'AUTH_MFA_LOCKED'
$ws9.Range("F7").Value = 'S2'
$ws9.Range("A8").Value = '감사 무결성 점검 FAIL'
$ws9.Range("B8").Value = 'chain hash 불연속'
$ws9.Range("C8").Value = 'chain-verify 결과 확인'
$ws9.Range("D8").Value = 'runbook_audit_chain 절차 실행'
$ws9.Range("E8").Value = '-'
$ws9.Range("F8").Value = 'S2'

# --- Sheet 10: 09_출처추적_매트릭스 ---
$ws10 = $wb.Worksheets.Item(10)
$ws10.Range("A2:F7").Style = "Normal"
$ws10.Range("A2").Value = 'M-001'
$ws10.Range("B2").Value = 'stale permission은 401'
$ws10.Range("C2").Value = 'docs/review/plans/20260221_auth_rbac_ops_admin_design_and_hardening_plan.md'
$ws10.Range("D2").Value = 'AUTH_STALE_PERMISSION'
$ws10.Range("E2").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F2").Value = 'PASS'
$ws10.Range("A3").Value = 'M-002'
$ws10.Range("B3").Value = 'lockout은 429'
$ws10.Range("C3").Value = 'docs/review/plans/20260221_auth_rbac_ops_admin_design_and_hardening_plan.md'
$ws10.Range("D3").Value = 'AUTH_LOCKED'
$ws10.Range("E3").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F3").Value = 'PASS'
$ws10.Range("A4").Value = 'M-003'
$ws10.Range("B4").Value = 'refresh reuse는 409'
$ws10.Range("C4").Value = 'docs/review/plans/20260221_auth_rbac_ops_admin_design_and_hardening_plan.md'
$ws10.Range("D4").Value = 'AUTH_REFRESH_REUSE_DETECTED'
$ws10.Range("E4").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F4").Value = 'PASS'
$ws10.Range("A5").Value = 'M-004'
$ws10.Range("B5").Value = 'PUBLIC/AUTHENTICATED는 access_level로 분리'
$ws10.Range("C5").Value = 'docs/review/plans/202603XX_go_live_gap_closure_plan.md'
$ws10.Range("D5").Value = 'access_level'
$ws10.Range("E5").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F5").Value = 'PASS'
$ws10.Range("A6").Value = 'M-005'
$ws10.Range("B6").Value = 'spec_consistency_check FAIL=0'
$ws10.Range("C6").Value = 'docs/review/mvp_verification_pack/artifacts/golive_spec_consistency_after.txt'
$ws10.Range("D6").Value = 'PASS=9 FAIL=0'
$ws10.Range("E6").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F6").Value = 'PASS'
$ws10.Range("A7").Value = 'M-006'
$ws10.Range("B7").Value = 'Node 22 표준화(.nvmrc)'
$ws10.Range("C7").Value = 'docs/reports/PROJECT_FULL_IMPLEMENTATION_AND_HARDENING_REPORT_202603XX.md'
$ws10.Range("D7").Value = '.nvmrc=22'
$ws10.Range("E7").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F7").Value = 'PASS'
$ws10.Range("A8").Value = 'M-007'
$ws10.Range("B8").Value = 'Audit Chain Verify API 반영'
$ws10.Range("C8").Value = 'docs/reports/PROJECT_FULL_IMPLEMENTATION_AND_HARDENING_REPORT_202603XX.md'
$ws10.Range("D8").Value = '/v1/admin/audit-logs/chain-verify'
$ws10.Range("E8").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F8").Value = 'PASS'
$ws10.Range("A9").Value = 'M-008'
$ws10.Range("B9").Value = 'Notion 동기화 완료 상태'
$ws10.Range("C9").Value = 'docs/review/mvp_verification_pack/artifacts/golive_notion_sync_status.txt'
$ws10.Range("D9").Value = 'Status: DONE'
$ws10.Range("E9").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F9").Value = 'PASS'
$ws10.Range("A10").Value = 'M-009'
$ws10.Range("B10").Value = 'Go-Live UTF-8 검증 통과'
$ws10.Range("C10").Value = 'docs/review/mvp_verification_pack/artifacts/golive_utf8_check.txt'
$ws10.Range("D10").Value = 'result=PASS'
$ws10.Range("E10").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F10").Value = 'PASS'
$ws10.Range("A11").Value = 'M-010'
$ws10.Range("B11").Value = '운영 runbook 3종 반영'
$ws10.Range("C11").Value = 'docs/ops/runbook_scheduler_lock.md'
$ws10.Range("D11").Value = 'runbook'
$ws10.Range("E11").Value = '2026-02-21T12:19:00Z'
$ws10.Range("F11").Value = 'PASS'

# --- Sheet 11: 10_문서변경이력 ---
$ws11 = $wb.Worksheets.Item(11)
$ws11.Range("A3").Value = 'v2026.02.21-golive'
$ws11.Range("B3").Value = '2026-02-21 21:17:31 +0900'
$ws11.Range("C3").Value = 'Go-Live 기준 반영(노드22, access_level 정규화, MFA/세션, RBAC 2인승인, audit chain verify, runbook, Notion 상태)'
$ws11.Range("D3").Value = 'S1,S2,S3,S4,S5,S6,S7'
$ws11.Range("A4").Value = 'v2026.02.21-golive'
$ws11.Range("B4").Value = '2026-02-21 21:19:00 +0900'
$ws11.Range("C4").Value = 'Go-Live 기준 반영(노드22, access_level 정규화, MFA/세션, RBAC 2인승인, audit chain verify, runbook, Notion 상태)'
$ws11.Range("D4").Value = 'S1,S2,S3,S4,S5,S6,S7'

